$d = $word.ActiveDocument

# --- Step 1: split paragraph before "it is harder" so our target run's
#     paragraph ends right after it, isolating it from the runs that
#     must stay untouched.
$full = $d.Content.Text
$idx = $full.IndexOf("it is harder")
$d.Range($idx, $idx).InsertParagraphBefore()

# --- Step 2: split before "squash" ---
$full = $d.Content.Text
$idx2 = $full.IndexOf("squash' to reduce")
$d.Range($idx2, $idx2).InsertParagraphBefore()

# --- Step 3: split right after "squash'" (word + following apostrophe) ---
$full = $d.Content.Text
$idx3 = $full.IndexOf("squash' to reduce") + 7
$d.Range($idx3, $idx3).InsertParagraphBefore()

# --- Step 4: replace the now-isolated "squash'" text with "rebase’" ---
$full = $d.Content.Text
$idx4 = $full.IndexOf("squash'")
$r = $d.Range($idx4, $idx4 + 7)
$r.Text = "rebase’"

# --- Step 5: merge the temporary paragraph breaks back, right to left,
#     preserving run structure on both sides of each merge point.
$full = $d.Content.Text
$idx5 = $full.IndexOf("it is harder")
$d.Range($idx5 - 1, $idx5).Delete()

$full = $d.Content.Text
$idx6 = $full.IndexOf("rebase’") + 7
$d.Range($idx6, $idx6 + 1).Delete()

$full = $d.Content.Text
$idx7 = $full.IndexOf("rebase’")
$d.Range($idx7 - 1, $idx7).Delete()

Write-Output "Done"
